# Auto update Excel log
# Appends newly logged sensor readings to the PIR, Humidity and Temperature
# sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cellRange, $text) {
    # Force the cell to keep its literal text representation instead of
    # letting Excel auto-convert date / percentage looking strings into
    # numeric values.
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

# ----------------------------------------------------------------------
# PIR sheet: rows 29-42 (Date, Timestamp, Hour, Location, Value, Status)
# ----------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("19:54:52", "No Motion",       "Inactive"),
    @("19:54:53", "No Motion",       "Inactive"),
    @("19:54:55", "No Motion",       "Inactive"),
    @("19:55:00", "No Motion",       "Inactive"),
    @("19:55:01", "Motion Detected", "Active"),
    @("19:55:08", "No Motion",       "Inactive"),
    @("19:55:13", "No Motion",       "Inactive"),
    @("19:55:18", "No Motion",       "Inactive"),
    @("19:55:24", "No Motion",       "Inactive"),
    @("19:55:29", "No Motion",       "Inactive"),
    @("19:55:34", "No Motion",       "Inactive"),
    @("19:55:38", "Motion Detected", "Active"),
    @("19:55:46", "No Motion",       "Inactive"),
    @("19:55:51", "No Motion",       "Inactive")
)

$startRow = 29
for ($i = 0; $i -lt $pirRows.Length; $i++) {
    $r = $startRow + $i
    $row = $pirRows[$i]

    Set-TextCell $wsPIR.Cells.Item($r, 1) "2026-02-01"
    $wsPIR.Cells.Item($r, 2).Value = $row[0]
    $wsPIR.Cells.Item($r, 3).Value = "19:00"
    $wsPIR.Cells.Item($r, 4).Value = "Bathroom"
    $wsPIR.Cells.Item($r, 5).Value = $row[1]
    $wsPIR.Cells.Item($r, 6).Value = $row[2]
}

# ----------------------------------------------------------------------
# Humidity sheet: rows 24-34 (Date, Timestamp, Hour, Location, Value, Status)
# ----------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("19:54:52", "78.1%"),
    @("19:54:52", "77.6%"),
    @("19:54:53", "78.7%"),
    @("19:54:55", "77.5%"),
    @("19:55:00", "78.7%"),
    @("19:55:10", "77.2%"),
    @("19:55:21", "77.2%"),
    @("19:55:31", "77.7%"),
    @("19:55:41", "78.6%"),
    @("19:55:46", "78.6%"),
    @("19:55:51", "78.6%")
)

$startRow = 24
for ($i = 0; $i -lt $humidityRows.Length; $i++) {
    $r = $startRow + $i
    $row = $humidityRows[$i]

    Set-TextCell $wsHumidity.Cells.Item($r, 1) "2026-02-01"
    $wsHumidity.Cells.Item($r, 2).Value = $row[0]
    $wsHumidity.Cells.Item($r, 3).Value = "19:00"
    $wsHumidity.Cells.Item($r, 4).Value = "Bathroom"
    Set-TextCell $wsHumidity.Cells.Item($r, 5) $row[1]
    $wsHumidity.Cells.Item($r, 6).Value = "Active"
}

# ----------------------------------------------------------------------
# Temperature sheet: rows 24-34 (Date, Timestamp, Hour, Location, Value, Status)
# ----------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")

$temperatureRows = @(
    @("19:54:52", "25.2C"),
    @("19:54:52", "25.2C"),
    @("19:54:53", "25.2C"),
    @("19:54:56", "25.2C"),
    @("19:55:01", "25.2C"),
    @("19:55:11", "25.2C"),
    @("19:55:21", "25.2C"),
    @("19:55:31", "25.2C"),
    @("19:55:41", "25.2C"),
    @("19:55:46", "25.2C"),
    @("19:55:51", "25.2C")
)

$startRow = 24
for ($i = 0; $i -lt $temperatureRows.Length; $i++) {
    $r = $startRow + $i
    $row = $temperatureRows[$i]

    Set-TextCell $wsTemperature.Cells.Item($r, 1) "2026-02-01"
    $wsTemperature.Cells.Item($r, 2).Value = $row[0]
    $wsTemperature.Cells.Item($r, 3).Value = "19:00"
    $wsTemperature.Cells.Item($r, 4).Value = "Bathroom"
    $wsTemperature.Cells.Item($r, 5).Value = $row[1]
    $wsTemperature.Cells.Item($r, 6).Value = "Active"
}
